$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17 and 18 have their "species record" data swapped (A,B,D,E,F,G,H,K,L,M,N,Q,R),
# while the shared observation metadata columns (C,I,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY)
# stay put.

# Capture row 17 values before overwriting (use Value2 - Value's getter is
# unreliable in this host and returns a reflection description string)
$A17 = $ws.Range("A17").Value2
$B17 = $ws.Range("B17").Value2
$D17 = $ws.Range("D17").Value2
$E17 = $ws.Range("E17").Value2
$F17 = $ws.Range("F17").Value2
$G17 = $ws.Range("G17").Value2
$H17 = $ws.Range("H17").Value2
$M17 = $ws.Range("M17").Value2
$Q17 = $ws.Range("Q17").Value2
$R17 = $ws.Range("R17").Value2

# Capture row 18 values before overwriting
$A18 = $ws.Range("A18").Value2
$B18 = $ws.Range("B18").Value2
$D18 = $ws.Range("D18").Value2
$E18 = $ws.Range("E18").Value2
$F18 = $ws.Range("F18").Value2
$G18 = $ws.Range("G18").Value2
$H18 = $ws.Range("H18").Value2
$Q18 = $ws.Range("Q18").Value2
$R18 = $ws.Range("R18").Value2

# Write row 18's former data into row 17
$ws.Range("A17").Value2 = $A18
$ws.Range("B17").Value2 = $B18
$ws.Range("D17").Value2 = $D18
$ws.Range("E17").Value2 = $E18
$ws.Range("F17").Value2 = $F18
$ws.Range("G17").Value2 = $G18
$ws.Range("H17").Value2 = $H18
$ws.Range("Q17").Value2 = $Q18
$ws.Range("R17").Value2 = $R18
# row17 loses its K/L/M/N content (row 18 had none there)
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

# Write row 17's former data into row 18
$ws.Range("A18").Value2 = $A17
$ws.Range("B18").Value2 = $B17
$ws.Range("D18").Value2 = $D17
$ws.Range("E18").Value2 = $E17
$ws.Range("F18").Value2 = $F17
$ws.Range("G18").Value2 = $G17
$ws.Range("H18").Value2 = $H17
$ws.Range("Q18").Value2 = $Q17
$ws.Range("R18").Value2 = $R17
# row18 gains the K/L/M/N content that used to be on row17 (K/L/N were blank
# placeholder cells there - touch NumberFormat to materialize an empty cell
# without writing literal text; M carried the activity text)
$ws.Range("K18").NumberFormat = "General"
$ws.Range("L18").NumberFormat = "General"
$ws.Range("M18").Value2 = $M17
$ws.Range("N18").NumberFormat = "General"
